$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (C) column dates from 2023-10-05 (45204) to 2023-10-08 (45207)
# for rows 2 through 8.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
